# Fruta / hortaliza, semanal
# Insert a new week's worth of data (2 rows: "Primera" and "Segunda" quality)
# at the top of the weekly series (rows 136-137), pushing all subsequent
# rows (old 136-163) down by two rows (new 138-165).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 136:163 down by inserting two blank rows at 136:137.
$ws.Rows("136:137").Insert()

# New row 136 (Calidad = Primera)
$ws.Range("A136").Value = 11
$ws.Range("B136").Value = "Vega Monumental Concepción"
$ws.Range("C136").Value = "Bíobío"
$ws.Range("D136").Value = 44798
$ws.Range("E136").Value = 8
$ws.Range("F136").Value = 100112044
$ws.Range("G136").Value = "Perejil"
$ws.Range("H136").Value = "Sin especificar"
$ws.Range("I136").Value = "Primera"
$ws.Range("J136").Value = 200
$ws.Range("K136").Value = 700
$ws.Range("L136").Value = 800
$ws.Range("M136").Value = 750
$ws.Range("N136").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O136").Value = "Región de Ñuble"
$ws.Range("P136").Value = 750
$ws.Range("Q136").Value = 1
$ws.Range("R136").Value = "Hortaliza"

# New row 137 (Calidad = Segunda)
$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 44798
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112044
$ws.Range("G137").Value = "Perejil"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Segunda"
$ws.Range("J137").Value = 100
$ws.Range("K137").Value = 600
$ws.Range("L137").Value = 600
$ws.Range("M137").Value = 600
$ws.Range("N137").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O137").Value = "Región de Ñuble"
$ws.Range("P137").Value = 600
$ws.Range("Q137").Value = 1
$ws.Range("R137").Value = "Hortaliza"
